# Page actions method added
# Update the "DATA" sheet's amazonHamburgerMenuTest rows so that the
# browser column uses "chrome" for rows 8 and 9 (previously "edge" and
# "firefox" respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()

$ws.Range("C8").Value = "chrome"
$ws.Range("C9").Value = "chrome"

$ws.Range("C7").Select()
